$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# --- Row 19: IUserRepository row gets marked done on 3/6/2019 ---
$ws.Range("C19").Value = "done"
$ws.Range("D19").Value = 43530
$ws.Range("D18").Copy()
$ws.Range("D19").PasteSpecial(-4122)

# --- Row 22: Logout method row now has status "b" (in Bearbeitung) ---
# Done before row 20 so the new "b" shared string is interned ahead of the
# corrected "Login Methode..." string (matches original authoring order).
$ws.Range("B22").Value = "Logout Methode im UserController erstellen"
$ws.Range("C22").Value = "b"

# --- Row 20: fix UserController typo, mark done on 3/6/2019 ---
$ws.Range("B20").Value = "Login Methode im UserController erstellen (inkl. Sessions)"
$ws.Range("C20").Value = "done"
$ws.Range("D20").Value = 43530
$ws.Range("D18").Copy()
$ws.Range("D20").PasteSpecial(-4122)

# --- Row 21: Login View row marked done on 3/6/2019 ---
$ws.Range("B21").Value = "Login View erstellen"
$ws.Range("C21").Value = "done"
$ws.Range("D21").Value = 43530
$ws.Range("D18").Copy()
$ws.Range("D21").PasteSpecial(-4122)

# --- Widen column B to fit the longer task descriptions ---
$ws.Range("B:B").ColumnWidth = 60.15

# --- Update the saved selection/scroll position to the bottom of the list ---
$ws.Activate()
[void]$ws.Range("D22").Select()
